$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Report")
$ws2 = $wb.Worksheets.Item("Quantities")

$ws1.Range("D12").Value = 300
$ws1.Range("E12").Value = 350
$ws1.Range("F12").Value = 310

$ws1.Range("D26").Value = 235

$ws1.Range("D27").Value = 80
$ws1.Range("E27").Value = 100
$ws1.Range("F27").Value = 90

$ws2.Range("C10").Value = 300
$ws2.Range("D10").Value = 320
$ws2.Range("E10").Value = 310

$ws2.Range("C11").Value = 234

$ws2.Range("F34").Select()
$ws1.Select()
$ws1.Range("E13").Select()
